$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1,
    1,
    3,
    1,
    2,
    3,
    3,
    1,
    2,
    1,
    2,
    3,
    3,
    1,
    1,
    3,
    3,
    2,
    3,
    1,
    3,
    1,
    3,
    1,
    1,
    2,
    1,
    1,
    3,
    2,
    3,
    3,
    1,
    3,
    1,
    1,
    3,
    3,
    1,
    3,
    3,
    3,
    2,
    3,
    2,
    3,
    1,
    2,
    2,
    2,
    3,
    3,
    3,
    2,
    1,
    3,
    3,
    2,
    2,
    3,
    2,
    3,
    1,
    1,
    3,
    3,
    1,
    3,
    3,
    2,
    3,
    3,
    2,
    3,
    2,
    3,
    1,
    2,
    3,
    3,
    1,
    2,
    1,
    3,
    3,
    2,
    1,
    3,
    1,
    3,
    3,
    1,
    2,
    1,
    3,
    3,
    3,
    3,
    3,
    3,
    3,
    2,
    3,
    3,
    3,
    1,
    3,
    3,
    2,
    3,
    2,
    1,
    1,
    1,
    3,
    3,
    1,
    3,
    2,
    1,
    3,
    3,
    1,
    3,
    1,
    2,
    3,
    3,
    2,
    3,
    1,
    3,
    3,
    3,
    2,
    1,
    3,
    1,
    2,
    3,
    3,
    3,
    2,
    3,
    1,
    3,
    2,
    3,
    3,
    3,
    1,
    1,
    3,
    1,
    3,
    1,
    1,
    3,
    1,
    1,
    1,
    1,
    1,
    3,
    2,
    1,
    2,
    3,
    1,
    3,
    2,
    3,
    3,
    3,
    1,
    1,
    3,
    3,
    3,
    1,
    2,
    3,
    2,
    3,
    1,
    3,
    1,
    1,
    3,
    1,
    3,
    1,
    3,
    3,
    1,
    3,
    1,
    1,
    2,
    1,
    2,
    3,
    2,
    1,
    1,
    3,
    3,
    3,
    1,
    3,
    3,
    2,
    3,
    3,
    1,
    3,
    3,
    3,
    1,
    3,
    1,
    1,
    3,
    3,
    3,
    2,
    3,
    2,
    1,
    3,
    3,
    2,
    1,
    2,
    3,
    3,
    3,
    1,
    2,
    2,
    1,
    2,
    1,
    1,
    2,
    1,
    3,
    2,
    2,
    2
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

